# Update the "取得日時" (retrieved datetime) timestamp for all existing data
# rows on the active sheet ("ランサーズ") from 2025-11-20 12:36:45 to
# 2025-11-20 12:48:22 (commit: "Append: 2025-11-20 12:48 JST").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2025-11-20 12:36:45"
$newTimestamp = "2025-11-20 12:48:22"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
